$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13
$ws.Cells.Item($row, 1).Value = 42620.889513888891
$ws.Cells.Item($row, 2).Value = 8
$ws.Cells.Item($row, 3).Value = 54
$ws.Cells.Item($row, 4).Value = 45
$ws.Cells.Item($row, 5).Value = 54
$ws.Cells.Item($row, 6).Value = 28
$ws.Cells.Item($row, 7).Value = 27915
$ws.Cells.Item($row, 8).Value = 13042
$ws.Cells.Item($row, 9).Value = 747
$ws.Cells.Item($row, 10).Value = 119
$ws.Cells.Item($row, 11).Value = 98
$ws.Cells.Item($row, 12).Value = 10
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Named"
